$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# Part 1: Append the two existing "New" sheet rows (old rows 2 and 3) onto the
# end of the "Previously added" sheet, as rows 137 and 138.
# ---------------------------------------------------------------------------

$sheet1Appended = @(
    @{ Row=137; A="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/indras-pag/kgfdd.html"; B="10 000 €";  C="Krāslava un raj."; D="3 ha.";  E="60620030054"; F=45929.39791666667 },
    @{ Row=138; A="https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/madona/dbnfi.html";        B="45 000 €";  C="Madona un raj.";   D="5 ha.";  E="70900080055"; F=45928.53611111111 }
)

foreach ($r in $sheet1Appended) {
    $row = $r.Row
    $ws1.Range("A$row").Value = $r.A
    $ws1.Range("B$row").Value = $r.B
    $ws1.Range("C$row").Value = $r.C
    $ws1.Range("D$row").Value = $r.D
    $ws1.Range("E$row").NumberFormat = "@"
    $ws1.Range("E$row").Value = $r.E
    $ws1.Range("F$row").Value = $r.F

    # Hyperlink first (it forces its own style), formatting copy afterwards
    $ws1.Hyperlinks.Add($ws1.Range("A$row"), $r.A) | Out-Null

    # Copy formatting (font/number-format/alignment) from the row above it,
    # which reuses the existing style indexes instead of creating new ones.
    $ws1.Range("A136:F136").Copy()
    $ws1.Range("A$row" + ":F$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Part 2: Replace the two rows on the "New" sheet with five freshly scraped
# listings (rows 2-6).
# ---------------------------------------------------------------------------

# Drop the hyperlinks that exist on the sheet today, so the newly added ones
# start again from rId1.
$ws2.Range("A2").Hyperlinks.Delete()

$sheet2New = @(
    @{ Row=2; A="https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/alsviku-pag/bdxkxj.html";     B="";          C="Alūksne un raj.";  D="";         E="";            F=45929.90277777778  },
    @{ Row=3; A="https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/sausnejas-pag/gkipe.html";     B="120 000 €"; C="Madona un raj.";   D="24 ha.";   E="70920010011"; F=45930.51458333334  },
    @{ Row=4; A="https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/cesvaines-l-t/adfpi.html";     B="59 000 €";  C="Madona un raj.";   D="13 ha.";   E="70270030008"; F=45930.50486111111  },
    @{ Row=5; A="https://www.ss.com/msg/lv/real-estate/wood/ogre-and-reg/mazozolu-pag/cenmm.html";        B="10 000 €";  C="Ogre un raj.";     D="3 ha.";    E="74720080031"; F=45930.429861111115 },
    @{ Row=6; A="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/lendzu-pag/ddghm.html";       B="8 100 €";   C="Rēzekne un raj.";  D="2.70 ha."; E="78660050212"; F=45930.433333333334 }
)

foreach ($r in $sheet2New) {
    $row = $r.Row
    $ws2.Range("A$row").Value = $r.A
    $ws2.Range("B$row").Value = $r.B
    $ws2.Range("C$row").Value = $r.C
    $ws2.Range("D$row").Value = $r.D
    $ws2.Range("E$row").NumberFormat = "@"
    $ws2.Range("E$row").Value = $r.E
    $ws2.Range("F$row").Value = $r.F

    $ws2.Hyperlinks.Add($ws2.Range("A$row"), $r.A) | Out-Null

    # Use the untouched row 136 of "Previously added" as the format template -
    # it carries exactly the same per-column styles (s=3/4/4/4/4/2) that every
    # data row on both sheets uses.
    $ws1.Range("A136:F136").Copy()
    $ws2.Range("A$row" + ":F$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Housekeeping: adding hyperlinks auto-registers a built-in "Hyperlink" cell
# style. All cells have since been re-painted with their original styles via
# PasteSpecial above, so nothing references it any more - drop it again.
# ---------------------------------------------------------------------------
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}
